{"js": "const body = context.document.body;\n\n// 1. Update the \"Test Limitations\" paragraph: detection-limit VAF changes\n//    from 2% (ASXL1 exception) to 4% (JAK2 exception), and the wording is\n//    consolidated into a single run.\nconst oldLimitations =\n  \"The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 2% with the exception of ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%-10%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. Synonymous variants are not routinely reported. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.\";\n\nconst newLimitations =\n  \"The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 4% with the exception of JAK2 c.1849G>T;p.(Val617Phe) (detection limit ~ 1%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. Synonymous variants are not routinely reported. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.\";\n\nconst limitationsResults = body.search(oldLimitations, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nlimitationsResults.load(\"items\");\nawait context.sync();\n\nif (limitationsResults.items.length > 0) {\n  limitationsResults.items[0].insertText(newLimitations, \"Replace\");\n  await context.sync();\n}\n\n// 2. Update the cached \"Reported\" date field result from 16-Sep-2024 to\n//    7-Oct-2024.\nconst dateResults = body.search(\"16-Sep-2024\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"7-Oct-2024\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the \"Test Limitations\" paragraph: detection-limit VAF changes\n#    from 2% (ASXL1 exception) to 4% (JAK2 exception), and the wording is\n#    consolidated into a single run.\n$oldLimitations = \"The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 2% with the exception of ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%-10%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. Synonymous variants are not routinely reported. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.\"\n\n$newLimitations = \"The detection limit of this assay for specimens sequenced to the target read depth of 250x is a variant allele frequency (VAF) of approximately 4% with the exception of JAK2 c.1849G>T;p.(Val617Phe) (detection limit ~ 1%). This assay is primarily qualitative however, the variant read frequency (VRF) is provided to assist with variant interpretation and is assumed to approximate VAF in most instances (noting that the VAF of some insertions/deletions may be underrepresented due to assay-based allele bias). Copy number variations, loss of heterozygosity, structural rearrangements or aneuploidies are not reported. Insertions or deletions (particularly those > 25 bp in length) are not reliably detected by this assay. Genes are analysed using the reference transcripts listed below; coding exons found in alternative transcripts are not assessed by this assay. This assay does not distinguish between somatic and germline variants. In addition, the clonal origin of somatic variants (i.e. disease compartment or cell lineage) cannot be determined. Synonymous variants are not routinely reported. Please note Peter Mac assumes sample identification, family relationships, and clinical diagnoses are as stated on the request. Our clinical recommendations may be based on evidence from third-party data sources and should be interpreted in the context of all other clinical and laboratory information for this patient.\"\n\n$d.Content.Find.Execute($oldLimitations, $false, $false, $false, $false, $false, $true, 1, $false, $newLimitations, 2)\n\n# 2. Update the cached \"Reported\" date field result from 16-Sep-2024 to\n#    7-Oct-2024.\n$d.Content.Find.Execute(\"16-Sep-2024\", $false, $false, $false, $false, $false, $true, 1, $false, \"7-Oct-2024\", 2)\n"}
